$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("stakeholder_a")
$wsA.Range("B5").Value = 0.3292956137742247
$wsA.Range("C5").Value = 0.3623351811588255
$wsA.Range("D5").Value = 0.3083692050669497
$wsA.Range("B8").Value = -0.197038665338371
$wsA.Range("C8").Value = 0.7923650869569785
$wsA.Range("B9").Value = -0.5876889617074239
$wsA.Range("C9").Value = -0.5668230331892988
$wsA.Range("B10").Value = 0.7847276270457945
$wsA.Range("C10").Value = -0.2255420537676797
$wsA.Range("B13").Value = 0.1003643930382884
$wsA.Range("C13").Value = 0.03913405214650325
$wsA.Range("B18").Value = 0.02336250214449371
$wsA.Range("C18").Value = 0.02105042454389415
$wsA.Range("D18").Value = 0.04853158610938179

$wsB = $wb.Worksheets.Item("stakeholder_b")
$wsB.Range("B5").Value = 0.3831879555466839
$wsB.Range("C5").Value = 0.2353965020895953
$wsB.Range("D5").Value = 0.3814155423637206
$wsB.Range("B8").Value = -0.7547402640566521
$wsB.Range("C8").Value = 0.3115024887193707
$wsB.Range("B9").Value = 0.107601063455275
$wsB.Range("C9").Value = -0.8093754862917217
$wsB.Range("B10").Value = 0.6471392006013769
$wsB.Range("C10").Value = 0.49787299757235
$wsB.Range("B13").Value = 0.2142149150689366
$wsB.Range("C13").Value = 0.08270015596735977
$wsB.Range("B18").Value = 0.09745855887588231
$wsB.Range("C18").Value = 0.05194144861041222
